$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" date column (C) for rows 2-9 from 45221 (2023-10-22)
# to 45224 (2023-10-25), keeping existing cell formatting.
for ($row = 2; $row -le 9; $row++) {
    $ws.Cells.Item($row, 3).Value = 45224
}
